{"js": "// Update the worksheet date and every \"NNN\u00d7N=\" multiplication prompt to the\n// new values from the latest generated output (commit 4250d90).\n//\n// Every original cell string is unique inside the document, so a simple\n// exact, case-sensitive search/replace per pair is unambiguous.\nconst replacements = [\n  [\"2024-06-15 Saturday\", \"2024-06-16 Sunday\"],\n  [\"407\u00d74=\", \"681\u00d79=\"],\n  [\"140\u00d73=\", \"171\u00d73=\"],\n  [\"746\u00d79=\", \"924\u00d77=\"],\n  [\"882\u00d74=\", \"981\u00d78=\"],\n  [\"769\u00d72=\", \"115\u00d76=\"],\n  [\"243\u00d77=\", \"301\u00d79=\"],\n  [\"520\u00d77=\", \"518\u00d77=\"],\n  [\"365\u00d78=\", \"333\u00d75=\"],\n  [\"653\u00d72=\", \"245\u00d78=\"],\n  [\"987\u00d78=\", \"464\u00d74=\"],\n  [\"159\u00d78=\", \"368\u00d73=\"],\n  [\"761\u00d73=\", \"582\u00d73=\"],\n  [\"690\u00d75=\", \"404\u00d75=\"],\n  [\"647\u00d72=\", \"313\u00d74=\"],\n  [\"139\u00d74=\", \"577\u00d74=\"],\n  [\"835\u00d74=\", \"709\u00d75=\"],\n  [\"545\u00d76=\", \"374\u00d77=\"],\n  [\"765\u00d79=\", \"450\u00d75=\"],\n  [\"969\u00d78=\", \"585\u00d72=\"],\n  [\"837\u00d79=\", \"608\u00d79=\"],\n  [\"911\u00d74=\", \"863\u00d76=\"],\n  [\"704\u00d72=\", \"863\u00d73=\"],\n  [\"385\u00d76=\", \"393\u00d76=\"],\n  [\"114\u00d78=\", \"858\u00d73=\"],\n  [\"125\u00d75=\", \"992\u00d73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every \"NNN\u00d7N=\" multiplication prompt to the\n# new values from the latest generated output (commit 4250d90).\n#\n# Every original cell string is unique inside the document, so a plain\n# Find/Replace (wdReplaceAll, but each string only ever matches once) per\n# pair is unambiguous and safe to run in any order.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-06-15 Saturday\", \"2024-06-16 Sunday\"),\n    @(\"407\u00d74=\", \"681\u00d79=\"),\n    @(\"140\u00d73=\", \"171\u00d73=\"),\n    @(\"746\u00d79=\", \"924\u00d77=\"),\n    @(\"882\u00d74=\", \"981\u00d78=\"),\n    @(\"769\u00d72=\", \"115\u00d76=\"),\n    @(\"243\u00d77=\", \"301\u00d79=\"),\n    @(\"520\u00d77=\", \"518\u00d77=\"),\n    @(\"365\u00d78=\", \"333\u00d75=\"),\n    @(\"653\u00d72=\", \"245\u00d78=\"),\n    @(\"987\u00d78=\", \"464\u00d74=\"),\n    @(\"159\u00d78=\", \"368\u00d73=\"),\n    @(\"761\u00d73=\", \"582\u00d73=\"),\n    @(\"690\u00d75=\", \"404\u00d75=\"),\n    @(\"647\u00d72=\", \"313\u00d74=\"),\n    @(\"139\u00d74=\", \"577\u00d74=\"),\n    @(\"835\u00d74=\", \"709\u00d75=\"),\n    @(\"545\u00d76=\", \"374\u00d77=\"),\n    @(\"765\u00d79=\", \"450\u00d75=\"),\n    @(\"969\u00d78=\", \"585\u00d72=\"),\n    @(\"837\u00d79=\", \"608\u00d79=\"),\n    @(\"911\u00d74=\", \"863\u00d76=\"),\n    @(\"704\u00d72=\", \"863\u00d73=\"),\n    @(\"385\u00d76=\", \"393\u00d76=\"),\n    @(\"114\u00d78=\", \"858\u00d73=\"),\n    @(\"125\u00d75=\", \"992\u00d73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
